$d = $word.ActiveDocument

# The "Class Blog instructions" block (a blank spacer paragraph, the
# "Class Blog instructions" heading, and its three ListParagraph bullet
# items) is being removed entirely. Locate it by content so the script is
# resilient to exact paragraph-index drift, then delete the run of
# paragraph marks + text from right after "Every social or policy
# problem..." through the end of the last bullet ("Respond to the
# question..."), leaving that first paragraph and the following
# (bookmark) paragraph directly adjacent.

$count = $d.Paragraphs.Count
$startPara = $null
$endPara = $null

for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "Every social or policy problem has a cost component*") {
        $startPara = $i
    }
    if ($t -like "Respond to the question by leaving a blog post*") {
        $endPara = $i
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    $rangeStart = $d.Paragraphs.Item($startPara).Range.End
    $rangeEnd = $d.Paragraphs.Item($endPara).Range.End

    $deleteRange = $d.Range($rangeStart, $rangeEnd)
    $deleteRange.Delete()
}

Write-Output "done"
